$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header cell A1 = "Category", formatted like the rest of row 1
# (bold/centered header style, same as B1:W1).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Category"

# The category cells in column A (rows 2..46) used to carry the header style
# (s="1"); the new layout keeps them as plain/default-styled cells.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 46
}
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Style = "Normal"
}
